# NIT-9008729253 Estado de Cuenta: remove previous EC data, add the new
# worker (HEBERT EDUARDO WATTS MATURANA) at the top, and replace Alberto
# Jose Arango Hernandez's period rows with the refreshed list (now sorted
# newest-period-first, 2005 -> 1706), per the updated database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First data row (16) is now the Hebert Eduardo Watts Maturana record.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143387239"
$ws.Range("D16").Value = "HEBERT EDUARDO WATTS MATURANA"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 9525
$ws.Range("G16").Value = 1236000

# Rows 17-52: Alberto Jose Arango Hernandez, one row per mora period, now
# listed newest-first (2005 down to 1706). Period 2005 keeps its special
# "Valor Mora" of 21143; every other period uses the standard 33384.
$albertoPeriods = @(
    "2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706"
)

$row = 17
foreach ($periodo in $albertoPeriods) {
    if ($periodo -eq "2005") {
        $valorMora = 21143
    } else {
        $valorMora = 33384
    }

    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "73133769"
    $ws.Range("D$row").Value = "ALBERTO JOSE ARANGO HERNANDEZ"
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = 834600

    $row = $row + 1
}
